$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7, shifting existing rows 7-15 down to 8-16.
$ws.Range("A7:R7").Insert()

# Populate the newly inserted row 7 with the new weekly price entry.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 45210
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112030
$ws.Range("G7").Value = "Poroto granado"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 550
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1600
$ws.Range("M7").Value = 1536
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 1536
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
